$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8935616612434387
$ws.Range("B1").Value = 1.649622797966003
$ws.Range("D1").Value = 1.552900433540344
$ws.Range("E1").Value = 1.015778541564941
